$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) to be treated as text so values such as
# "232.76" or "19.00" are stored verbatim instead of being converted
# into numbers (which would drop trailing zeros / round the value).
$ws.Range("D2").Value = "37.253.19"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "2.058.52"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.76"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  +1.90%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.46"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.381"
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.96"
$ws.Range("E10").Value = "  -0.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0759"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "2.363.08"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.57"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.58"
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.774"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.11"
$ws.Range("E17").Value = "  -2.50%  "
$ws.Range("D18").Value = "2.061.62"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").Value = "37.192.20"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.30"
$ws.Range("E20").Value = "  +6.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.19"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").Value = "0.0₃0808"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "225.79"
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.38"
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.81"
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("E28").Value = "  +5.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.74"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("E30").Value = "  -2.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.00"
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("E32").Value = "  -0.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.42"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0615"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.57"
$ws.Range("E35").Value = "  +4.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.49"
$ws.Range("E36").Value = "  -1.34%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.23"
$ws.Range("E39").Value = "  -2.46%  "
$ws.Range("E40").Value = "  -5.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.96"
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").Value = "1.468.28"
$ws.Range("E42").Value = "  -0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0937"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "95.73"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0212"
$ws.Range("E45").Value = "  +1.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.29"
$ws.Range("E46").Value = "  -1.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.16"
$ws.Range("E47").Value = "  +2.61%  "
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.04"
$ws.Range("E49").Value = "  -6.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.14"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("E51").Value = "  +0.32%  "
